$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 7155.075790473336
$ws.Range("O2").Value = 6980.325566461754

# Sheet "2030" (sheet2)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 5707.815717280662
$ws.Range("I2").Value = 44492.05901988943
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 21991.42050229464
$ws.Range("O2").Value = 12080.30369160272

# Sheet "2035" (sheet3)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15119.1038010506
$ws.Range("O2").Value = 14761.94879382113

# Sheet "2040" (sheet4)
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15226.22770907155
$ws.Range("O2").Value = 14761.94879382113

# Sheet "2045" (sheet5)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15770.95526052178
$ws.Range("O2").Value = 17097.41478016988

# Sheet "2050" (sheet6)
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15770.95526052178
$ws.Range("O2").Value = 17097.41478016988
